# CIERRE 7 SEPT 2022
# - Update the incentive-month label on "VALES DE INSENTIVOS" from JUNIO to AGOSTO
# - Make "VALES DE INSENTIVOS" the active/selected tab (was "ARQUITECTO")
# - TODAY() cells recalc automatically on save

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Update the month/year text for the incentive payment
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE  AGOSTO   2022"

# Switch the active sheet from "ARQUITECTO" to "VALES DE INSENTIVOS"
$wsVales.Activate()
